$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = -93.23
$ws.Range("D6").Value = 30.06
$ws.Range("J6").Value = 30.06
$ws.Range("D9").Value = 30.06
$ws.Range("J9").Value = 30.06
$ws.Range("B11").Value = -93.23
